$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency values (price/volume refresh, plus two row re-orderings)
$ws.Range('D2').Value = '26.092.65'
$ws.Range('E2').Value = '  +0.72%  '

$ws.Range('D3').Value = '1.647.28'
$ws.Range('E3').Value = '  +1.00%  '

$ws.Range('E4').Value = '  +1.00%  '

$ws.Range('D5').Value = '217.18'
$ws.Range('E5').Value = '  +1.26%  '

$ws.Range('E6').Value = '  +1.14%  '

$ws.Range('E7').Value = '  +1.03%  '

$ws.Range('D8').Value = '0.257'
$ws.Range('E8').Value = '  +0.63%  '

$ws.Range('D9').Value = '0.0641'
$ws.Range('E9').Value = '  +1.54%  '

$ws.Range('D10').Value = '19.70'
$ws.Range('E10').Value = '  +0.35%  '

$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').Value = '  +0.98%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '4.30'
$ws.Range('E12').Value = '  +1.66%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.873.19'
$ws.Range('E13').Value = '  +0.95%  '

$ws.Range('D14').Value = '1.626.05'
$ws.Range('E14').Value = '  -0.88%  '

$ws.Range('D15').Value = '0.546'
$ws.Range('E15').Value = '  +0.31%  '

$ws.Range('D16').Value = '0.0₃0766'
$ws.Range('E16').Value = '  +1.34%  '

$ws.Range('D17').Value = '63.47'
$ws.Range('E17').Value = '  +1.11%  '

$ws.Range('D18').Value = '26.168.56'
$ws.Range('E18').Value = '  +1.10%  '

$ws.Range('E19').Value = '  +1.09%  '

$ws.Range('D20').Value = '193.62'
$ws.Range('E20').Value = '  +0.29%  '

$ws.Range('D21').Value = '4.35'
$ws.Range('E21').Value = '  -0.75%  '

$ws.Range('D22').Value = '9.96'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').Value = '6.25'
$ws.Range('E23').Value = '  -0.32%  '

$ws.Range('D24').Value = '1.82'
$ws.Range('E24').Value = '  +1.58%  '

$ws.Range('D25').Value = '144.63'
$ws.Range('E25').Value = '  +1.71%  '

$ws.Range('E26').Value = '  +1.36%  '

$ws.Range('E27').Value = '  +4.12%  '

$ws.Range('D28').Value = '6.92'
$ws.Range('E28').Value = '  +0.92%  '

$ws.Range('D29').Value = '15.57'
$ws.Range('E29').Value = '  +0.80%  '

$ws.Range('E30').Value = '  +1.41%  '

$ws.Range('D31').Value = '0.0499'
$ws.Range('E31').Value = '  -0.25%  '

$ws.Range('E32').Value = '  -0.47%  '

$ws.Range('D33').Value = '3.28'
$ws.Range('E33').Value = '  +1.70%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '1.53'
$ws.Range('E34').Value = '  -2.94%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.47'
$ws.Range('E35').Value = '  +2.13%  '

$ws.Range('D36').Value = '0.908'
$ws.Range('E36').Value = '  +0.78%  '

$ws.Range('D37').Value = '1.136.32'

$ws.Range('D38').Value = '0.542'
$ws.Range('E38').Value = '  -1.49%  '

$ws.Range('E39').Value = '  +0.08%  '

$ws.Range('E40').Value = '  +0.68%  '

$ws.Range('D41').Value = '5.51'
$ws.Range('E41').Value = '  +1.14%  '

$ws.Range('D42').Value = '99.72'
$ws.Range('E42').Value = '  +0.52%  '

$ws.Range('D43').Value = '0.799'
$ws.Range('E43').Value = '  -0.62%  '

$ws.Range('D44').Value = '1.782.41'
$ws.Range('E44').Value = '  +1.01%  '

$ws.Range('E45').Value = '  +4.63%  '

$ws.Range('D46').Value = '56.77'
$ws.Range('E46').Value = '  +1.16%  '

$ws.Range('E47').Value = '  +0.98%  '

$ws.Range('D48').Value = '1.46'
$ws.Range('E48').Value = '  +0.05%  '

$ws.Range('D49').Value = '7.74'
$ws.Range('E49').Value = '  +1.55%  '

$ws.Range('E50').Value = '  +0.62%  '

$ws.Range('D51').Value = '0.0959'
$ws.Range('E51').Value = '  -0.24%  '
